# Auto update Excel log
# Appends new sensor/event rows to the Proximity, mmWave and Camera logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Proximity sheet: append rows 57-63 (Living Room Main Door ENTER/EXIT)
# ---------------------------------------------------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")

$proximityRows = @(
    @("2026-01-30", "14:31:02", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-01-30", "14:31:05", "14:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-01-30", "14:31:24", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-01-30", "14:31:48", "14:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-01-30", "14:32:13", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-01-30", "14:32:16", "14:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-01-30", "14:32:19", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
)

$startRow = 57
for ($i = 0; $i -lt $proximityRows.Count; $i++) {
    $r = $startRow + $i
    $row = $proximityRows[$i]
    # Prefix the Date column with a quote so Excel keeps it as literal text
    # instead of auto-converting it to a date serial number.
    $wsProximity.Cells.Item($r, 1).Value = "'" + $row[0]
    $wsProximity.Cells.Item($r, 2).Value = $row[1]
    $wsProximity.Cells.Item($r, 3).Value = $row[2]
    $wsProximity.Cells.Item($r, 4).Value = $row[3]
    $wsProximity.Cells.Item($r, 5).Value = $row[4]
    $wsProximity.Cells.Item($r, 6).Value = $row[5]
}

# ---------------------------------------------------------------------
# mmWave sheet: append rows 78-81 (Living Room FALL_DETECTED / EMERGENCY)
# ---------------------------------------------------------------------
$wsMmWave = $wb.Worksheets.Item("mmWave")

$mmWaveRows = @(
    @("2026-01-30", "14:30:59", "14:00", "Living Room", "FALL_DETECTED", "EMERGENCY"),
    @("2026-01-30", "14:30:59", "14:00", "Living Room", "FALL_DETECTED", "EMERGENCY"),
    @("2026-01-30", "14:31:23", "14:00", "Living Room", "FALL_DETECTED", "EMERGENCY"),
    @("2026-01-30", "14:31:23", "14:00", "Living Room", "FALL_DETECTED", "EMERGENCY")
)

$startRow = 78
for ($i = 0; $i -lt $mmWaveRows.Count; $i++) {
    $r = $startRow + $i
    $row = $mmWaveRows[$i]
    $wsMmWave.Cells.Item($r, 1).Value = "'" + $row[0]
    $wsMmWave.Cells.Item($r, 2).Value = $row[1]
    $wsMmWave.Cells.Item($r, 3).Value = $row[2]
    $wsMmWave.Cells.Item($r, 4).Value = $row[3]
    $wsMmWave.Cells.Item($r, 5).Value = $row[4]
    $wsMmWave.Cells.Item($r, 6).Value = $row[5]
}

# ---------------------------------------------------------------------
# Camera sheet: append rows 3-8 (Living Room Main Door Image Captured)
# ---------------------------------------------------------------------
$wsCamera = $wb.Worksheets.Item("Camera")

$cameraRows = @(
    @("2026-01-30", "14:31:02", "14:00", "Living Room Main Door", "Image Captured", "Event=ENTER"),
    @("2026-01-30", "14:31:25", "14:00", "Living Room Main Door", "Image Captured", "Event=ENTER"),
    @("2026-01-30", "14:31:48", "14:00", "Living Room Main Door", "Image Captured", "Event=ENTER"),
    @("2026-01-30", "14:32:13", "14:00", "Living Room Main Door", "Image Captured", "Event=ENTER"),
    @("2026-01-30", "14:32:15", "14:00", "Living Room Main Door", "Image Captured", "Event=ENTER"),
    @("2026-01-30", "14:32:20", "14:00", "Living Room Main Door", "Image Captured", "Event=ENTER")
)

$startRow = 3
for ($i = 0; $i -lt $cameraRows.Count; $i++) {
    $r = $startRow + $i
    $row = $cameraRows[$i]
    $wsCamera.Cells.Item($r, 1).Value = "'" + $row[0]
    $wsCamera.Cells.Item($r, 2).Value = $row[1]
    $wsCamera.Cells.Item($r, 3).Value = $row[2]
    $wsCamera.Cells.Item($r, 4).Value = $row[3]
    $wsCamera.Cells.Item($r, 5).Value = $row[4]
    $wsCamera.Cells.Item($r, 6).Value = $row[5]
}
